$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.423.09"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.560.49"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.03"
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.500"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.05"
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.249"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0591"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0865"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.780.22"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.565.76"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.81"
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.516"
$ws.Range("E15").Value = "  -2.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.44"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.401.26"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "212.64"
$ws.Range("E18").Value = "  -3.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0687"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.24"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.07"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.52"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("E24").Value = "  +3.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.02"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.67"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.95"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("E29").Value = "  -2.14%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.20"
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.371.08"
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.96"
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.53"
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.954"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.819"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.975"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.80"
$ws.Range("E43").Value = "  +1.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.69"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.23"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.694.13"
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.15"
$ws.Range("E48").Value = "  -3.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₇0995"
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0953"
$ws.Range("E51").Value = "  -1.91%  "
